$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before current row 52 (shifts existing rows 52-123 down to 53-124)
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new observation
$ws.Cells.Item(52, 1).Value = 11
$ws.Cells.Item(52, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(52, 3).Value = "Bíobío"
$ws.Cells.Item(52, 4).Value = 44546
$ws.Cells.Item(52, 5).Value = 8
$ws.Cells.Item(52, 6).Value = "Fruta"
$ws.Cells.Item(52, 7).Value = 100108
$ws.Cells.Item(52, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(52, 9).Value = 100108005
$ws.Cells.Item(52, 10).Value = "Piña"
$ws.Cells.Item(52, 11).Value = "Caramelo"
$ws.Cells.Item(52, 12).Value = "Primera"
$ws.Cells.Item(52, 13).Value = 100
$ws.Cells.Item(52, 14).Value = 17000
$ws.Cells.Item(52, 15).Value = 18000
$ws.Cells.Item(52, 16).Value = 17500
$ws.Cells.Item(52, 17).Value = "$/caja 12 unidades"
$ws.Cells.Item(52, 18).Value = "Ecuador"
$ws.Cells.Item(52, 19).Value = 1458
$ws.Cells.Item(52, 20).Value = 12
